$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.898.26"
$ws.Range("E2").Value = "  +2.20%  "

$ws.Range("D3").Value = "2.049.32"
$ws.Range("E3").Value = "  +1.37%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.38"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.615"
$ws.Range("E6").Value = "  +1.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "57.96"
$ws.Range("E7").Value = "  +5.55%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.385"
$ws.Range("E9").Value = "  +2.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0805"
$ws.Range("E10").Value = "  +2.60%  "

$ws.Range("E11").Value = "  +1.17%  "

$ws.Range("D12").Value = "2.352.14"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.52"
$ws.Range("E13").Value = "  +2.79%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.57"
$ws.Range("E14").Value = "  +1.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.751"
$ws.Range("E15").Value = "  +1.33%  "

$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("D17").Value = "2.046.30"
$ws.Range("E17").Value = "  +1.41%  "

$ws.Range("D18").Value = "37.842.65"
$ws.Range("E18").Value = "  +2.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.12"
$ws.Range("E19").Value = "  -1.37%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.58"
$ws.Range("E20").Value = "  +0.90%  "

$ws.Range("D21").Value = "0.0₃0827"
$ws.Range("E21").Value = "  +1.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.42"
$ws.Range("E22").Value = "  +0.40%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.45"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.24"
$ws.Range("E25").Value = "  +2.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "166.19"
$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.26"
$ws.Range("E27").Value = "  +0.80%  "

$ws.Range("E28").Value = "  +6.71%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.98"
$ws.Range("E29").Value = "  +1.23%  "

$ws.Range("E30").Value = "  -0.45%  "

$ws.Range("E31").Value = "  +1.50%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.51"
$ws.Range("E32").Value = "  +0.31%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.56"
$ws.Range("E33").Value = "  +3.01%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0610"
$ws.Range("E34").Value = "  -0.48%  "

$ws.Range("E35").Value = "  +10.29%  "

$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.95"
$ws.Range("E37").Value = "  +12.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +4.73%  "

$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("E40").Value = "  +0.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.57"
$ws.Range("E41").Value = "  +2.19%  "

$ws.Range("D42").Value = "1.482.70"
$ws.Range("E42").Value = "  +0.53%  "

$ws.Range("E43").Value = "  +3.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0932"
$ws.Range("E44").Value = "  +1.60%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.64"
$ws.Range("E45").Value = "  +2.14%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.25"
$ws.Range("E46").Value = "  +19.41%  "

$ws.Range("E47").Value = "  +0.07%  "

$ws.Range("E48").Value = "  +0.45%  "

$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.99"
$ws.Range("E50").Value = "  -2.80%  "

$ws.Range("D51").Value = "2.244.77"
$ws.Range("E51").Value = "  +1.85%  "
